$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1566.6666
$ws.Range("J43").Value = 1475
$ws.Range("L43").Value = 1475
$ws.Range("N43").Value = -1613
$ws.Range("H62").Value = 4668.3335
$ws.Range("I62").Value = 5302.1
$ws.Range("K62").Value = 5302.1
$ws.Range("M62").Value = -4678.1
$ws.Range("H65").Value = 4668.3335
$ws.Range("I65").Value = 5302.1
$ws.Range("K65").Value = 26510.5
$ws.Range("M65").Value = -23390.5
$ws.Range("H80").Value = 1564.3529
$ws.Range("I80").Value = 1571.4546
$ws.Range("J80").Value = 1551.3334
$ws.Range("K80").Value = 4714.3638
$ws.Range("L80").Value = 4654.0002
$ws.Range("M80").Value = -3716.3638
$ws.Range("N80").Value = -6650.0002
$ws.Range("H83").Value = 1564.3529
$ws.Range("I83").Value = 1571.4546
$ws.Range("J83").Value = 1551.3334
$ws.Range("K83").Value = 14143.0914
$ws.Range("L83").Value = 13962.0006
$ws.Range("M83").Value = -9151.091400000001
$ws.Range("N83").Value = -23946.0006
$ws.Range("H132").Value = 897.87177
$ws.Range("I132").Value = 849.75
$ws.Range("K132").Value = 2549.25
$ws.Range("M132").Value = -19.25
$ws.Range("H137").Value = 35529.55
$ws.Range("I137").Value = 764.8333
$ws.Range("J137").Value = 92417.27
$ws.Range("K137").Value = 2294.4999
$ws.Range("L137").Value = 277251.81
$ws.Range("M137").Value = 255.5001000000002
$ws.Range("N137").Value = -282351.81
$ws.Range("H139").Value = 69811.336
$ws.Range("J139").Value = 69811.336
$ws.Range("L139").Value = 69811.336
$ws.Range("N139").Value = -80091.336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 231822.42
$ws.Range("I2").Value = 278056.34
$ws.Range("J2").Value = 652.75
$ws.Range("K2").Value = 278056.34
$ws.Range("L2").Value = 652.75
$ws.Range("M2").Value = -277943.34
$ws.Range("N2").Value = -878.75
$ws.Range("H32").Value = 1731.95
$ws.Range("I32").Value = 1688.8384
$ws.Range("J32").Value = 6000
$ws.Range("K32").Value = 1688.8384
$ws.Range("L32").Value = 6000
$ws.Range("M32").Value = -1401.8384
$ws.Range("N32").Value = -6574
$ws.Range("H110").Value = 1164
$ws.Range("I110").Value = 1159.5
$ws.Range("K110").Value = 1159.5
$ws.Range("M110").Value = 885.5
$ws.Range("H116").Value = 231822.42
$ws.Range("I116").Value = 278056.34
$ws.Range("J116").Value = 652.75
$ws.Range("K116").Value = 278056.34
$ws.Range("L116").Value = 652.75
$ws.Range("M116").Value = -275762.34
$ws.Range("N116").Value = -5240.75
$ws.Range("H122").Value = 2630.5715
$ws.Range("I122").Value = 1402.3334
$ws.Range("K122").Value = 4207.0002
$ws.Range("M122").Value = -1757.0002
$ws.Range("H132").Value = 1768.2559
$ws.Range("I132").Value = 1383.4517
$ws.Range("K132").Value = 4150.355100000001
$ws.Range("M132").Value = -1620.355100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 231822.42
$ws.Range("I3").Value = 278056.34
$ws.Range("J3").Value = 652.75
$ws.Range("K3").Value = 278056.34
$ws.Range("L3").Value = 652.75
$ws.Range("M3").Value = -277942.34
$ws.Range("N3").Value = -880.75
$ws.Range("H105").Value = 2386.1333
$ws.Range("J105").Value = 2613.75
$ws.Range("L105").Value = 2613.75
$ws.Range("N105").Value = -6107.75
$ws.Range("H107").Value = 564.8333
$ws.Range("I107").Value = 461.54544
$ws.Range("J107").Value = 727.1429000000001
$ws.Range("K107").Value = 461.54544
$ws.Range("L107").Value = 727.1429000000001
$ws.Range("M107").Value = 1458.45456
$ws.Range("N107").Value = -4567.1429
$ws.Range("H140").Value = 47114.145
$ws.Range("J140").Value = 47114.145
$ws.Range("L140").Value = 47114.145
$ws.Range("N140").Value = -57474.145

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 917
$ws.Range("I16").Value = 917
$ws.Range("K16").Value = 917
$ws.Range("M16").Value = -630
$ws.Range("H31").Value = 2458.4211
$ws.Range("J31").Value = 3355.875
$ws.Range("L31").Value = 3355.875
$ws.Range("N31").Value = -3945.875
$ws.Range("H34").Value = 2458.4211
$ws.Range("J34").Value = 3355.875
$ws.Range("L34").Value = 3355.875
$ws.Range("N34").Value = -3759.875
$ws.Range("H105").Value = 1009.6667
$ws.Range("I105").Value = 1012.9167
$ws.Range("J105").Value = 996.6667
$ws.Range("K105").Value = 1012.9167
$ws.Range("L105").Value = 996.6667
$ws.Range("M105").Value = 734.0833
$ws.Range("N105").Value = -4490.6667
$ws.Range("H107").Value = 959.86664
$ws.Range("I107").Value = 1127.9166
$ws.Range("J107").Value = 287.66666
$ws.Range("K107").Value = 1127.9166
$ws.Range("L107").Value = 287.66666
$ws.Range("M107").Value = 792.0834
$ws.Range("N107").Value = -4127.66666
$ws.Range("H113").Value = 917
$ws.Range("I113").Value = 917
$ws.Range("K113").Value = 917
$ws.Range("M113").Value = 1253
$ws.Range("H122").Value = 2543.2856
$ws.Range("I122").Value = 1464.6666
$ws.Range("J122").Value = 4484.8
$ws.Range("K122").Value = 4393.9998
$ws.Range("L122").Value = 13454.4
$ws.Range("M122").Value = -1943.9998
$ws.Range("N122").Value = -18354.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1416639.2
$ws.Range("I4").Value = 1708468.4
$ws.Range("J4").Value = 249323
$ws.Range("K4").Value = 5125405.199999999
$ws.Range("L4").Value = 747969
$ws.Range("M4").Value = -5125293.199999999
$ws.Range("N4").Value = -748193
$ws.Range("H5").Value = 505.8
$ws.Range("I5").Value = 472.3
$ws.Range("K5").Value = 1416.9
$ws.Range("M5").Value = -1304.9
$ws.Range("H23").Value = 99.125
$ws.Range("J23").Value = 138
$ws.Range("L23").Value = 414
$ws.Range("N23").Value = -884
$ws.Range("H131").Value = 17077.465
$ws.Range("J131").Value = 18020.906
$ws.Range("L131").Value = 54062.71799999999
$ws.Range("N131").Value = -64142.71799999999
$ws.Range("H133").Value = 3908.25
$ws.Range("I133").Value = 2374.75
$ws.Range("J133").Value = 4675
$ws.Range("K133").Value = 7124.25
$ws.Range("L133").Value = 14025
$ws.Range("M133").Value = -2064.25
$ws.Range("N133").Value = -24145
$ws.Range("H135").Value = 505.8
$ws.Range("I135").Value = 472.3
$ws.Range("K135").Value = 4250.7
$ws.Range("M135").Value = -1715.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1307.2
$ws.Range("I122").Value = 1159.9
$ws.Range("K122").Value = 3479.7
$ws.Range("M122").Value = -1029.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3637
$ws.Range("I61").Value = 3368
$ws.Range("J61").Value = 4444
$ws.Range("K61").Value = 3368
$ws.Range("L61").Value = 4444
$ws.Range("M61").Value = -3166
$ws.Range("N61").Value = -4848
$ws.Range("H113").Value = 3637
$ws.Range("I113").Value = 3368
$ws.Range("J113").Value = 4444
$ws.Range("K113").Value = 3368
$ws.Range("L113").Value = 4444
$ws.Range("M113").Value = -1198
$ws.Range("N113").Value = -8784
$ws.Range("H132").Value = 1941.1163
$ws.Range("I132").Value = 1348.5
$ws.Range("K132").Value = 4045.5
$ws.Range("M132").Value = -1515.5
$ws.Range("H136").Value = 3894.25
$ws.Range("I136").Value = 3215.125
$ws.Range("K136").Value = 9645.375
$ws.Range("M136").Value = -7095.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 807.7
$ws.Range("I107").Value = 675.2222
$ws.Range("K107").Value = 2025.6666
$ws.Range("M107").Value = -105.6666
$ws.Range("H122").Value = 30043.822
$ws.Range("I122").Value = 39396.047
$ws.Range("K122").Value = 118188.141
$ws.Range("M122").Value = -115738.141
$ws.Range("H132").Value = 1354.1025
$ws.Range("I132").Value = 1150.9412
$ws.Range("K132").Value = 3452.8236
$ws.Range("M132").Value = -922.8235999999997
